# Auto-generated edit script: updates column F (interested-count) values
# across the four worksheets, per the source diff (regenerated site data / gh-pages output).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 35
$ws.Range("F4").Value = 8369
$ws.Range("F5").Value = 8369
$ws.Range("F8").Value = 531
$ws.Range("F9").Value = 7340
$ws.Range("F10").Value = 589
$ws.Range("F16").Value = 170
$ws.Range("F17").Value = 123
$ws.Range("F18").Value = 148
$ws.Range("F19").Value = 12174
$ws.Range("F22").Value = 2468
$ws.Range("F23").Value = 3536
$ws.Range("F24").Value = 53
$ws.Range("F26").Value = 2915
$ws.Range("F27").Value = 110
$ws.Range("F30").Value = 3352
$ws.Range("F33").Value = 1718
$ws.Range("F36").Value = 6017
$ws.Range("F38").Value = 1832
$ws.Range("F39").Value = 1258
$ws.Range("F41").Value = 902
$ws.Range("F45").Value = 198
$ws.Range("F48").Value = 1582

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 28
$ws.Range("F15").Value = 8
$ws.Range("F20").Value = 919
$ws.Range("F22").Value = 74

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 321
$ws.Range("F3").Value = 468

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 35
$ws.Range("F4").Value = 28
$ws.Range("F5").Value = 321
$ws.Range("F8").Value = 8369
$ws.Range("F11").Value = 531
$ws.Range("F12").Value = 7340
$ws.Range("F13").Value = 7340
$ws.Range("F14").Value = 589
$ws.Range("F20").Value = 170
$ws.Range("F21").Value = 148
$ws.Range("F23").Value = 12174
$ws.Range("F27").Value = 2468
$ws.Range("F28").Value = 2468
$ws.Range("F29").Value = 3536
$ws.Range("F30").Value = 110
$ws.Range("F33").Value = 8
$ws.Range("F34").Value = 3352
$ws.Range("F37").Value = 1718
$ws.Range("F40").Value = 6017
$ws.Range("F41").Value = 74
$ws.Range("F42").Value = 1832
$ws.Range("F44").Value = 1258
$ws.Range("F46").Value = 902
$ws.Range("F48").Value = 198
$ws.Range("F51").Value = 1582
